$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update projected_kcp (column D) values for rows 205-366 (season weeks 30-52)
# per corrected cultivar projection data.

for ($r = 205; $r -le 211; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4540283
}

for ($r = 212; $r -le 218; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4558204
}

for ($r = 219; $r -le 225; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4571071
}

for ($r = 226; $r -le 232; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4578865
}

for ($r = 233; $r -le 239; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4578663
}

for ($r = 240; $r -le 246; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4565129
}

for ($r = 247; $r -le 253; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4531566
}

for ($r = 254; $r -le 260; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4471211
}

for ($r = 261; $r -le 267; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4378224
}

for ($r = 268; $r -le 274; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4248062
}

for ($r = 275; $r -le 281; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.4077434
}

for ($r = 282; $r -le 288; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.3864316
}

for ($r = 289; $r -le 295; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.3608309
}

for ($r = 296; $r -le 302; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.3311315
}

for ($r = 303; $r -le 309; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.2978617
}

for ($r = 310; $r -le 316; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.2620732
}

for ($r = 317; $r -le 323; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.2255623
}

for ($r = 324; $r -le 330; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.1908385
}

for ($r = 331; $r -le 337; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.1605056
}

for ($r = 338; $r -le 344; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.1362931
}

for ($r = 345; $r -le 351; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.1185185
}

for ($r = 352; $r -le 358; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.1063507
}

for ($r = 359; $r -le 366; $r++) {
    $ws.Cells.Item($r, 4).Value = 0.0980685
}
